$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet (canonical URL build date + concept count) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-03-28T11:28:20+02:00"
$meta.Range("B22").Value = "'206"

# --- Append new concept rows to the Concepts sheet ---
$ws = $wb.Worksheets.Item("Concepts")

# Copy the formatting of an existing data row onto the new row range, then fill in values
$fmtSrc = $ws.Range("A2:D2")
$fmtSrc.Copy()
$ws.Range("A107:D207").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A107").Value = "'1"
$ws.Range("B107").Value = 'O13-OTH'
$ws.Range("C107").Value = 'Other'
$ws.Range("D107").Value = 'Other'
$ws.Range("A108").Value = "'1"
$ws.Range("B108").Value = 'O7-1'
$ws.Range("C108").Value = 'Affymetrix Cytoscan HD'
$ws.Range("A109").Value = "'1"
$ws.Range("B109").Value = 'O7-2'
$ws.Range("C109").Value = 'Affymetrix EMET Plus Premier Pack'
$ws.Range("A110").Value = "'1"
$ws.Range("B110").Value = 'O7-3'
$ws.Range("C110").Value = 'Affymetrix Genome-Wide Human SNP Array 5.0'
$ws.Range("A111").Value = "'1"
$ws.Range("B111").Value = 'O7-4'
$ws.Range("C111").Value = 'Affymetrix Genome-Wide Human SNP Array 6.0'
$ws.Range("A112").Value = "'1"
$ws.Range("B112").Value = 'O7-5'
$ws.Range("C112").Value = 'Affymetrix HT Human Genome U133A Array Plate Set'
$ws.Range("A113").Value = "'1"
$ws.Range("B113").Value = 'O7-6'
$ws.Range("C113").Value = 'Affymetrix Human Exon 1.0 ST'
$ws.Range("A114").Value = "'1"
$ws.Range("B114").Value = 'O7-7'
$ws.Range("C114").Value = 'Affymetrix Human Gene 1.0 ST'
$ws.Range("A115").Value = "'1"
$ws.Range("B115").Value = 'O7-8'
$ws.Range("C115").Value = 'Affymetrix Human Genome U219 Array Plate'
$ws.Range("A116").Value = "'1"
$ws.Range("B116").Value = 'O7-9'
$ws.Range("C116").Value = 'Affymetrix Human MIP 330K'
$ws.Range("A117").Value = "'1"
$ws.Range("B117").Value = 'O7-10'
$ws.Range("C117").Value = 'Affymetrix Human U133 Plus 2.0'
$ws.Range("A118").Value = "'1"
$ws.Range("B118").Value = 'O7-11'
$ws.Range("C118").Value = 'Affymetrix Human U133 Plus PM'
$ws.Range("A119").Value = "'1"
$ws.Range("B119").Value = 'O7-12'
$ws.Range("C119").Value = 'Affymetrix Mapping 100K Array Set'
$ws.Range("A120").Value = "'1"
$ws.Range("B120").Value = 'O7-13'
$ws.Range("C120").Value = 'Affymetrix Mapping 10K 2.0 Array Set'
$ws.Range("A121").Value = "'1"
$ws.Range("B121").Value = 'O7-14'
$ws.Range("C121").Value = 'Affymetrix Mapping 500K Array Set'
$ws.Range("A122").Value = "'1"
$ws.Range("B122").Value = 'O7-15'
$ws.Range("C122").Value = 'Affymetrix OncoScan FFPE Express 2.0'
$ws.Range("A123").Value = "'1"
$ws.Range("B123").Value = 'O7-16'
$ws.Range("C123").Value = 'Agilent 244K Custom Gene Expression G4502A-07'
$ws.Range("A124").Value = "'1"
$ws.Range("B124").Value = 'O7-17'
$ws.Range("C124").Value = 'Agilent 244K Custom Gene Expression G4502A-07-1'
$ws.Range("A125").Value = "'1"
$ws.Range("B125").Value = 'O7-18'
$ws.Range("C125").Value = 'Agilent 244K Custom Gene Expression G4502A-07-2'
$ws.Range("A126").Value = "'1"
$ws.Range("B126").Value = 'O7-19'
$ws.Range("C126").Value = 'Agilent 244K Custom Gene Expression G4502A-07-3'
$ws.Range("A127").Value = "'1"
$ws.Range("B127").Value = 'O7-20'
$ws.Range("C127").Value = 'Agilent 8 x 15K Human miRNA-specific microarray'
$ws.Range("A128").Value = "'1"
$ws.Range("B128").Value = 'O7-21'
$ws.Range("C128").Value = 'Agilent Human CGH 1x1M'
$ws.Range("A129").Value = "'1"
$ws.Range("B129").Value = 'O7-22'
$ws.Range("C129").Value = 'Agilent Human CGH 2x400K'
$ws.Range("A130").Value = "'1"
$ws.Range("B130").Value = 'O7-23'
$ws.Range("C130").Value = 'Agilent Human CGH 4x180K'
$ws.Range("A131").Value = "'1"
$ws.Range("B131").Value = 'O7-24'
$ws.Range("C131").Value = 'Agilent Human CGH 8x60K'
$ws.Range("A132").Value = "'1"
$ws.Range("B132").Value = 'O7-25'
$ws.Range("C132").Value = 'Agilent Human CNV 2x400K'
$ws.Range("A133").Value = "'1"
$ws.Range("B133").Value = 'O7-26'
$ws.Range("C133").Value = 'Agilent Human CNV Association 2x105K'
$ws.Range("A134").Value = "'1"
$ws.Range("B134").Value = 'O7-27'
$ws.Range("C134").Value = 'Agilent Human CpG Island Microarray Kit'
$ws.Range("A135").Value = "'1"
$ws.Range("B135").Value = 'O7-28'
$ws.Range("C135").Value = 'Agilent Human Genome 105A'
$ws.Range("A136").Value = "'1"
$ws.Range("B136").Value = 'O7-29'
$ws.Range("C136").Value = 'Agilent Human Genome 244A'
$ws.Range("A137").Value = "'1"
$ws.Range("B137").Value = 'O7-30'
$ws.Range("C137").Value = 'Agilent Human Genome 44K'
$ws.Range("A138").Value = "'1"
$ws.Range("B138").Value = 'O7-31'
$ws.Range("C138").Value = 'Agilent Human Genome CGH Custom Microaary 2x415K'
$ws.Range("A139").Value = "'1"
$ws.Range("B139").Value = 'O7-32'
$ws.Range("C139").Value = 'Agilent Human miRNA Microarray Kit (v2)'
$ws.Range("A140").Value = "'1"
$ws.Range("B140").Value = 'O7-33'
$ws.Range("C140").Value = 'Agilent Human Promoter ChIP-on-chip Microarray Set'
$ws.Range("A141").Value = "'1"
$ws.Range("B141").Value = 'O7-34'
$ws.Range("C141").Value = 'Agilent Human SpliceArray'
$ws.Range("A142").Value = "'1"
$ws.Range("B142").Value = 'O7-35'
$ws.Range("C142").Value = 'Agilent Whole Human Genome Oligo Microarray Kit'
$ws.Range("A143").Value = "'1"
$ws.Range("B143").Value = 'O7-36'
$ws.Range("C143").Value = 'Almac Human CRC'
$ws.Range("A144").Value = "'1"
$ws.Range("B144").Value = 'O7-37'
$ws.Range("C144").Value = 'capillary sequencing'
$ws.Range("A145").Value = "'1"
$ws.Range("B145").Value = 'O7-38'
$ws.Range("C145").Value = 'Complete Genomics'
$ws.Range("A146").Value = "'1"
$ws.Range("B146").Value = 'O7-39'
$ws.Range("C146").Value = 'Custom-designed cDNA array'
$ws.Range("A147").Value = "'1"
$ws.Range("B147").Value = 'O7-40'
$ws.Range("C147").Value = 'Custom-designed gene expression array'
$ws.Range("A148").Value = "'1"
$ws.Range("B148").Value = 'O7-41'
$ws.Range("C148").Value = 'Helicos sequencing'
$ws.Range("A149").Value = "'1"
$ws.Range("B149").Value = 'O7-42'
$ws.Range("C149").Value = 'HumanOmni2.5-8 BeadChip Kit'
$ws.Range("A150").Value = "'1"
$ws.Range("B150").Value = 'O7-43'
$ws.Range("C150").Value = 'Illumina GA sequencing'
$ws.Range("A151").Value = "'1"
$ws.Range("B151").Value = 'O7-44'
$ws.Range("C151").Value = 'Illumina goldengate methylation'
$ws.Range("A152").Value = "'1"
$ws.Range("B152").Value = 'O7-45'
$ws.Range("C152").Value = 'Illumina GoldenGate Methylation Cancer Panel I'
$ws.Range("A153").Value = "'1"
$ws.Range("B153").Value = 'O7-46'
$ws.Range("C153").Value = 'Illumina HiSeq'
$ws.Range("A154").Value = "'1"
$ws.Range("B154").Value = 'O7-47'
$ws.Range("C154").Value = 'Illumina HiSeq X Ten'
$ws.Range("A155").Value = "'1"
$ws.Range("B155").Value = 'O7-48'
$ws.Range("C155").Value = 'Illumina Human Omni1-Quad beadchip'
$ws.Range("A156").Value = "'1"
$ws.Range("B156").Value = 'O7-49'
$ws.Range("C156").Value = 'Illumina Human1M OmniQuad chip'
$ws.Range("A157").Value = "'1"
$ws.Range("B157").Value = 'O7-50'
$ws.Range("C157").Value = 'Illumina human1m-duo'
$ws.Range("A158").Value = "'1"
$ws.Range("B158").Value = 'O7-51'
$ws.Range("C158").Value = 'Illumina human510s-duo'
$ws.Range("A159").Value = "'1"
$ws.Range("B159").Value = 'O7-52'
$ws.Range("C159").Value = 'Illumina human660w-quad'
$ws.Range("A160").Value = "'1"
$ws.Range("B160").Value = 'O7-53'
$ws.Range("C160").Value = 'Illumina HumanCNV370-Duo v1.0 BeadChip'
$ws.Range("A161").Value = "'1"
$ws.Range("B161").Value = 'O7-54'
$ws.Range("C161").Value = 'Illumina humancytosnp-12'
$ws.Range("A162").Value = "'1"
$ws.Range("B162").Value = 'O7-55'
$ws.Range("C162").Value = 'Illumina HumanHap550'
$ws.Range("A163").Value = "'1"
$ws.Range("B163").Value = 'O7-56'
$ws.Range("C163").Value = 'Illumina HumanHT-12 v4.0 beadchip'
$ws.Range("A164").Value = "'1"
$ws.Range("B164").Value = 'O7-57'
$ws.Range("C164").Value = 'Illumina humanht-16'
$ws.Range("A165").Value = "'1"
$ws.Range("B165").Value = 'O7-58'
$ws.Range("C165").Value = 'Illumina humanht-17'
$ws.Range("A166").Value = "'1"
$ws.Range("B166").Value = 'O7-59'
$ws.Range("C166").Value = 'Illumina humanmethylation27'
$ws.Range("A167").Value = "'1"
$ws.Range("B167").Value = 'O7-60'
$ws.Range("C167").Value = 'Illumina HumanOmniExpress BeadChip'
$ws.Range("A168").Value = "'1"
$ws.Range("B168").Value = 'O7-61'
$ws.Range("C168").Value = 'Illumina HumanRef-8 v3.0 beadchip'
$ws.Range("A169").Value = "'1"
$ws.Range("B169").Value = 'O7-62'
$ws.Range("C169").Value = 'Illumina HumanWG-6 v3.0 beadchip'
$ws.Range("A170").Value = "'1"
$ws.Range("B170").Value = 'O7-63'
$ws.Range("C170").Value = 'Illumina Infinium HumanMethylation450'
$ws.Range("A171").Value = "'1"
$ws.Range("B171").Value = 'O7-64'
$ws.Range("C171").Value = 'Illumina microRNA Expression Profiling Panel'
$ws.Range("A172").Value = "'1"
$ws.Range("B172").Value = 'O7-65'
$ws.Range("C172").Value = 'Illumina MiSeq Personal Sequencer'
$ws.Range("A173").Value = "'1"
$ws.Range("B173").Value = 'O7-66'
$ws.Range("C173").Value = 'Ion Torrent PGM'
$ws.Range("A174").Value = "'1"
$ws.Range("B174").Value = 'O7-67'
$ws.Range("C174").Value = 'Ion Torrent Proton'
$ws.Range("A175").Value = "'1"
$ws.Range("B175").Value = 'O7-68'
$ws.Range("C175").Value = 'M.D. Anderson Reverse Phase Protein Array Core'
$ws.Range("A176").Value = "'1"
$ws.Range("B176").Value = 'O7-69'
$ws.Range("C176").Value = 'Microsatellite Instability Analysis'
$ws.Range("A177").Value = "'1"
$ws.Range("B177").Value = 'O7-70'
$ws.Range("C177").Value = 'nanoString'
$ws.Range("A178").Value = "'1"
$ws.Range("B178").Value = 'O7-71'
$ws.Range("C178").Value = 'Nimblegen CGS'
$ws.Range("A179").Value = "'1"
$ws.Range("B179").Value = 'O7-72'
$ws.Range("C179").Value = 'Nimblegen Gene Expression 12x135K'
$ws.Range("A180").Value = "'1"
$ws.Range("B180").Value = 'O7-73'
$ws.Range("C180").Value = 'Nimblegen Gene Expression 385K'
$ws.Range("A181").Value = "'1"
$ws.Range("B181").Value = 'O7-74'
$ws.Range("C181").Value = 'Nimblegen Gene Expression 4x72K'
$ws.Range("A182").Value = "'1"
$ws.Range("B182").Value = 'O7-75'
$ws.Range("C182").Value = 'Nimblegen Human CGH 2.1M Whole-Genome v2.0D Array'
$ws.Range("A183").Value = "'1"
$ws.Range("B183").Value = 'O7-76'
$ws.Range("C183").Value = 'Nimblegen Human CGH 3x720 Whole-Genome v3.0 Array'
$ws.Range("A184").Value = "'1"
$ws.Range("B184").Value = 'O7-77'
$ws.Range("C184").Value = 'Nimblegen Human Methylation 2.1M Whole-Genome sets'
$ws.Range("A185").Value = "'1"
$ws.Range("B185").Value = 'O7-78'
$ws.Range("C185").Value = 'Nimblegen Human Methylation 385K Whole-Genome sets'
$ws.Range("A186").Value = "'1"
$ws.Range("B186").Value = 'O7-79'
$ws.Range("C186").Value = 'PacBio RS sequencing'
$ws.Range("A187").Value = "'1"
$ws.Range("B187").Value = 'O7-80'
$ws.Range("C187").Value = 'PCR'
$ws.Range("A188").Value = "'1"
$ws.Range("B188").Value = 'O7-81'
$ws.Range("C188").Value = 'PCR and capillary sequencing'
$ws.Range("A189").Value = "'1"
$ws.Range("B189").Value = 'O7-82'
$ws.Range("C189").Value = 'qPCR'
$ws.Range("A190").Value = "'1"
$ws.Range("B190").Value = 'O7-83'
$ws.Range("C190").Value = 'Roche 454 sequencing'
$ws.Range("A191").Value = "'1"
$ws.Range("B191").Value = 'O7-84'
$ws.Range("C191").Value = 'Sequenom MassARRAY'
$ws.Range("A192").Value = "'1"
$ws.Range("B192").Value = 'O7-85'
$ws.Range("C192").Value = 'SOLiD sequencing'
$ws.Range("A193").Value = "'1"
$ws.Range("B193").Value = 'O7-86'
$ws.Range("C193").Value = 'Digital PCR'
$ws.Range("A194").Value = "'1"
$ws.Range("B194").Value = 'O7-87'
$ws.Range("C194").Value = 'NextSeq 500'
$ws.Range("A195").Value = "'1"
$ws.Range("B195").Value = 'O7-88'
$ws.Range("C195").Value = 'NovaSeq'
$ws.Range("A196").Value = "'1"
$ws.Range("B196").Value = 'O7-OTH'
$ws.Range("C196").Value = 'Other'
$ws.Range("A197").Value = "'1"
$ws.Range("B197").Value = 'O10-1'
$ws.Range("C197").Value = 'Ion AmpliSeq Cancer Hotspot Panel v2'
$ws.Range("A198").Value = "'1"
$ws.Range("B198").Value = 'O10-2'
$ws.Range("C198").Value = 'Ion AmpliSeq Colon and Lung Cancer Research Panel v2'
$ws.Range("A199").Value = "'1"
$ws.Range("B199").Value = 'O10-3'
$ws.Range("C199").Value = 'Ion AmpliSeq Comprehensive Cancer Panel'
$ws.Range("A200").Value = "'1"
$ws.Range("B200").Value = 'O10-4'
$ws.Range("C200").Value = 'Ion AmpliSeq Oncomine Comprehensive Assay'
$ws.Range("A201").Value = "'1"
$ws.Range("B201").Value = 'O10-5'
$ws.Range("C201").Value = 'Ion AmpliSeq Oncomine Focus Assay'
$ws.Range("A202").Value = "'1"
$ws.Range("B202").Value = 'O10-6'
$ws.Range("C202").Value = 'Ion AmpliSeq TP53 Research Panel'
$ws.Range("A203").Value = "'1"
$ws.Range("B203").Value = 'O10-7'
$ws.Range("C203").Value = 'Lyric: Ion AmpliSeq Profiler'
$ws.Range("A204").Value = "'1"
$ws.Range("B204").Value = 'O10-8'
$ws.Range("C204").Value = 'QIAGEN GeneRead DNAseq Targeted Panels V2'
$ws.Range("A205").Value = "'1"
$ws.Range("B205").Value = 'O10-9'
$ws.Range("C205").Value = 'SAFIR02 Panel'
$ws.Range("A206").Value = "'1"
$ws.Range("B206").Value = 'O10-10'
$ws.Range("C206").Value = 'Mosc3'
$ws.Range("A207").Value = "'1"
$ws.Range("B207").Value = '010-OTH'
$ws.Range("C207").Value = 'Other'
